$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ C = 0.5379009246826172;  D = 67.27800000000001 }
    3  = @{ C = 0.16324782371521;    D = 38.991 }
    4  = @{ C = 1.078111886978149;   D = 56.173 }
    5  = @{ C = 0.3983108997344971;  D = 73.245 }
    6  = @{ C = 0.1716568470001221;  D = 6.017 }
    7  = @{ C = 0.171022891998291;   D = 6.017 }
    8  = @{ C = 0.4400262832641602;  D = 54.218 }
    9  = @{ C = 0.3820259571075439;  D = 71.29000000000001 }
    10 = @{ C = 0.171309232711792;   D = 7.972 }
    11 = @{ C = 0.1744680404663086;  D = 7.972 }
    12 = @{ C = 0.3911838531494141;  D = 56.173 }
    13 = @{ C = 0.3939950466156006;  D = 73.245 }
    14 = @{ C = 0.1747791767120361;  D = 64.398 }
    15 = @{ C = 0.1765859127044678;  D = 522.806 }
    16 = @{ C = 0.5177819728851318;  D = 58.093 }
    17 = @{ C = 0.7613980770111084;  D = 522.668 }
    18 = @{ C = 0.184607982635498;   D = 6.215 }
    19 = @{ C = 0.184988260269165;   D = 6.215 }
    20 = @{ C = 0.4199428558349609;  D = 56.138 }
    21 = @{ C = 0.7898678779602051;  D = 520.713 }
    22 = @{ C = 0.1952550411224365;  D = 8.17 }
    23 = @{ C = 0.1950888633728027;  D = 8.17 }
    24 = @{ C = 0.4299178123474121;  D = 58.093 }
    25 = @{ C = 0.5953450202941895;  D = 522.668 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("D$row").Value = $values[$row].D
}
